$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3125
$ws.Range("J40").Value = 3500
$ws.Range("L40").Value = 3500
$ws.Range("N40").Value = -3850
$ws.Range("H53").Value = 1362.5
$ws.Range("I53").Value = 2618.75
$ws.Range("K53").Value = 2618.75
$ws.Range("M53").Value = -1981.75
$ws.Range("H63").Value = 18271
$ws.Range("J63").Value = 18271
$ws.Range("L63").Value = 18271
$ws.Range("N63").Value = -19519
$ws.Range("H66").Value = 18271
$ws.Range("J66").Value = 18271
$ws.Range("L66").Value = 54813
$ws.Range("N66").Value = -61053
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("H86").Value = 1290.0667
$ws.Range("I86").Value = 1239.3572
$ws.Range("K86").Value = 1239.3572
$ws.Range("M86").Value = -116.3571999999999
$ws.Range("H88").Value = 1173.2858
$ws.Range("J88").Value = 1481.5
$ws.Range("L88").Value = 1481.5
$ws.Range("N88").Value = -2293.5
$ws.Range("H89").Value = 1290.0667
$ws.Range("I89").Value = 1239.3572
$ws.Range("K89").Value = 6196.786
$ws.Range("M89").Value = -580.7860000000001
$ws.Range("H91").Value = 1173.2858
$ws.Range("J91").Value = 1481.5
$ws.Range("L91").Value = 1481.5
$ws.Range("N91").Value = -4289.5
$ws.Range("H96").Value = 1235.5454
$ws.Range("I96").Value = 962.75
$ws.Range("J96").Value = 1963
$ws.Range("K96").Value = 2888.25
$ws.Range("L96").Value = 5889
$ws.Range("M96").Value = -1515.25
$ws.Range("N96").Value = -8635
$ws.Range("H107").Value = 2076.3333
$ws.Range("I107").Value = 1339.4
$ws.Range("K107").Value = 1339.4
$ws.Range("M107").Value = 580.5999999999999
$ws.Range("H137").Value = 46695.863
$ws.Range("I137").Value = 866.5
$ws.Range("K137").Value = 2599.5
$ws.Range("M137").Value = -49.5
$ws.Range("H141").Value = 4669772.5
$ws.Range("I141").Value = 9334064
$ws.Range("J141").Value = 5481.3335
$ws.Range("K141").Value = 28002192
$ws.Range("L141").Value = 16444.0005
$ws.Range("M141").Value = -27997012
$ws.Range("N141").Value = -26804.0005

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5926.392
$ws.Range("I32").Value = 3569.9343
$ws.Range("K32").Value = 3569.9343
$ws.Range("M32").Value = -3282.9343
$ws.Range("H61").Value = 53908.438
$ws.Range("I61").Value = 67652.914
$ws.Range("J61").Value = 12675
$ws.Range("K61").Value = 67652.914
$ws.Range("L61").Value = 12675
$ws.Range("M61").Value = -67440.914
$ws.Range("N61").Value = -13099
$ws.Range("H74").Value = 834.5484
$ws.Range("I74").Value = 587.3333
$ws.Range("K74").Value = 587.3333
$ws.Range("M74").Value = 286.6667
$ws.Range("H77").Value = 834.5484
$ws.Range("I77").Value = 587.3333
$ws.Range("K77").Value = 2936.6665
$ws.Range("M77").Value = 1431.3335
$ws.Range("H132").Value = 2409.6155
$ws.Range("I132").Value = 3262.8333
$ws.Range("J132").Value = 2153.65
$ws.Range("K132").Value = 9788.499899999999
$ws.Range("L132").Value = 6460.950000000001
$ws.Range("M132").Value = -7258.499899999999
$ws.Range("N132").Value = -11520.95
$ws.Range("H136").Value = 53908.438
$ws.Range("I136").Value = 67652.914
$ws.Range("J136").Value = 12675
$ws.Range("K136").Value = 202958.742
$ws.Range("L136").Value = 38025
$ws.Range("M136").Value = -200408.742
$ws.Range("N136").Value = -43125

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2891.2703
$ws.Range("I20").Value = 3078.7273
$ws.Range("J20").Value = 2616.3333
$ws.Range("K20").Value = 3078.7273
$ws.Range("L20").Value = 2616.3333
$ws.Range("M20").Value = -2831.7273
$ws.Range("N20").Value = -3110.3333
$ws.Range("H105").Value = 2713.1
$ws.Range("I105").Value = 2489.8333
$ws.Range("K105").Value = 2489.8333
$ws.Range("M105").Value = -742.8332999999998
$ws.Range("H132").Value = 50375
$ws.Range("J132").Value = 50375
$ws.Range("L132").Value = 50375
$ws.Range("N132").Value = -60495
$ws.Range("H134").Value = 7237.2964
$ws.Range("I134").Value = 9135.888999999999
$ws.Range("J134").Value = 3440.111
$ws.Range("K134").Value = 27407.667
$ws.Range("L134").Value = 10320.333
$ws.Range("M134").Value = -24872.667
$ws.Range("N134").Value = -15390.333

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 368
$ws.Range("I7").Value = 526
$ws.Range("K7").Value = 526
$ws.Range("M7").Value = -413
$ws.Range("H58").Value = 1403913.2
$ws.Range("I58").Value = 3346191.5
$ws.Range("J58").Value = 1156.8334
$ws.Range("K58").Value = 3346191.5
$ws.Range("L58").Value = 1156.8334
$ws.Range("M58").Value = -3345988.5
$ws.Range("N58").Value = -1562.8334
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("H105").Value = 1201.125
$ws.Range("I105").Value = 1214.5333
$ws.Range("K105").Value = 1214.5333
$ws.Range("M105").Value = 532.4666999999999
$ws.Range("H122").Value = 3483.353
$ws.Range("I122").Value = 1936.875
$ws.Range("K122").Value = 5810.625
$ws.Range("M122").Value = -3360.625
$ws.Range("H134").Value = 2973.8667
$ws.Range("I134").Value = 2926.125
$ws.Range("J134").Value = 3028.4285
$ws.Range("K134").Value = 8778.375
$ws.Range("L134").Value = 9085.2855
$ws.Range("M134").Value = -6243.375
$ws.Range("N134").Value = -14155.2855
$ws.Range("H136").Value = 1403913.2
$ws.Range("I136").Value = 3346191.5
$ws.Range("J136").Value = 1156.8334
$ws.Range("K136").Value = 10038574.5
$ws.Range("L136").Value = 3470.5002
$ws.Range("M136").Value = -10036024.5
$ws.Range("N136").Value = -8570.5002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1382656.9
$ws.Range("I4").Value = 1577074.4
$ws.Range("K4").Value = 4731223.199999999
$ws.Range("M4").Value = -4731111.199999999
$ws.Range("H5").Value = 552.0454999999999
$ws.Range("J5").Value = 774.75
$ws.Range("L5").Value = 2324.25
$ws.Range("N5").Value = -2548.25
$ws.Range("H113").Value = 38921.242
$ws.Range("J113").Value = 1067.9131
$ws.Range("L113").Value = 3203.7393
$ws.Range("N113").Value = -7543.7393
$ws.Range("H135").Value = 552.0454999999999
$ws.Range("J135").Value = 774.75
$ws.Range("L135").Value = 6972.75
$ws.Range("N135").Value = -12042.75

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492
$ws.Range("H110").Value = 77601
$ws.Range("J110").Value = 77601
$ws.Range("L110").Value = 77601
$ws.Range("N110").Value = -85781
$ws.Range("H113").Value = 1176.75
$ws.Range("I113").Value = 799
$ws.Range("K113").Value = 799
$ws.Range("M113").Value = 1371
$ws.Range("H132").Value = 3850557.5
$ws.Range("I132").Value = 7696512
$ws.Range("J132").Value = 4603
$ws.Range("K132").Value = 23089536
$ws.Range("L132").Value = 13809
$ws.Range("M132").Value = -23087006
$ws.Range("N132").Value = -18869

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1003.7917
$ws.Range("I22").Value = 489
$ws.Range("K22").Value = 489
$ws.Range("M22").Value = -194
$ws.Range("H27").Value = 1003.7917
$ws.Range("I27").Value = 489
$ws.Range("K27").Value = 489
$ws.Range("M27").Value = -382
$ws.Range("H46").Value = 2923.3635
$ws.Range("I46").Value = 1800
$ws.Range("K46").Value = 1800
$ws.Range("M46").Value = -1612
$ws.Range("H100").Value = 1897.7142
$ws.Range("I100").Value = 957
$ws.Range("K100").Value = 957
$ws.Range("M100").Value = -416
$ws.Range("H106").Value = 30000
$ws.Range("J106").Value = 30000
$ws.Range("L106").Value = 30000
$ws.Range("N106").Value = -32524
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 2093.2559
$ws.Range("I132").Value = 1748.591
$ws.Range("J132").Value = 2454.3333
$ws.Range("K132").Value = 5245.772999999999
$ws.Range("L132").Value = 7362.999899999999
$ws.Range("M132").Value = -2715.772999999999
$ws.Range("N132").Value = -12422.9999
$ws.Range("H136").Value = 2572.7646
$ws.Range("I136").Value = 1811.0834
$ws.Range("J136").Value = 4400.8
$ws.Range("K136").Value = 5433.2502
$ws.Range("L136").Value = 13202.4
$ws.Range("M136").Value = -2883.2502
$ws.Range("N136").Value = -18302.4

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 28866
$ws.Range("J16").Value = 28866
$ws.Range("L16").Value = 28866
$ws.Range("N16").Value = -29450
$ws.Range("H96").Value = 1706.4736
$ws.Range("I96").Value = 1700.6666
$ws.Range("K96").Value = 1700.6666
$ws.Range("M96").Value = -327.6666
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H113").Value = 1049.5714
$ws.Range("I113").Value = 965.1667
$ws.Range("K113").Value = 2895.5001
$ws.Range("M113").Value = -725.5001000000002
$ws.Range("H136").Value = 37038940
$ws.Range("I136").Value = 50506580
$ws.Range("J136").Value = 2933.75
$ws.Range("K136").Value = 151519740
$ws.Range("L136").Value = 8801.25
$ws.Range("M136").Value = -151517190
$ws.Range("N136").Value = -13901.25
